$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 804, shifting existing rows 804:845 down to 805:846
$ws.Rows.Item(804).Insert()

# Write the new row's values. Column A holds a date formatted as plain text
# (e.g. "2026/02/13"), so force text entry to avoid Excel's automatic
# date-literal conversion, then restore the default "Normal" style so no
# stray number-format style gets attached to the cell.
$ws.Cells.Item(804, 1).NumberFormat = "@"
$ws.Cells.Item(804, 1).Value = "2026/02/13"
$ws.Cells.Item(804, 1).Style = "Normal"

$ws.Cells.Item(804, 2).Value = "金"
$ws.Cells.Item(804, 3).Value = 16
$ws.Cells.Item(804, 4).Value = 201
